$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, centered, bordered) from an existing header
# cell onto the three new header cells before we stamp in the labels.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record (Wins / Losses / Ties) for every player row.
for ($row = 2; $row -le 39; $row++) {
    $ws.Cells.Item($row, 30).Value = 105
    $ws.Cells.Item($row, 31).Value = 57
    $ws.Cells.Item($row, 32).Value = 0
}

Write-Output "Updated season record columns AD:AF for rows 1-39"
